$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new loc row at 45 -------------------------------------------
# Insert() copies the formatting of the row above (row 44: styles 5/6/6),
# which matches the styling the new row should have.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).RowHeight = 46.5

$ws.Cells.Item(45,1).Value = "cwl_warn_deserialize_ele_final"
$ws.Cells.Item(45,3).Value = "this is the last time CWL warns about it, the rest will be silently ignored to reduce log spam"
$ws.Cells.Item(45,4).Value = "this is the last time CWL warns about it, the rest will be silently ignored to reduce log spam"

# --- Update the safety-cone description (now row 47 after the insert) ----
$newSafetyDesc = "This element is missing or modified from your current game. `nCWL kept the game going by replacing it with a safety cone.`nYou should report the relevant information to mod author or CWL.`nUsing this element will let CWL purge it from your save.`nYou may also keep the safety cone, CWL will restore it when the responsible mod functions again."
$ws.Cells.Item(47,3).Value = $newSafetyDesc
$ws.Cells.Item(47,4).Value = $newSafetyDesc
$ws.Rows.Item(47).RowHeight = 186

# --- Restore view / selection state ---------------------------------------
[void]$ws.Range("A41").Select()
$excel.ActiveWindow.ScrollRow = 41
[void]$ws.Range("A45").Select()
